$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intents")

# 1. Append the new keywords to the AskAboutRestaurant "Text" cell (B13)
$oldText = $ws.Range("B13").Value2
$ws.Range("B13").Value = $oldText + ",email.website,phone,telephone"

# 2. Widen column B so the expanded text fits (stored width ends up at 201
#    once the engine's character->pixel->character round trip is applied).
$ws.Columns.Item(2).ColumnWidth = 200.28571428571428

# 3. Update the active selection / scroll position to reflect the edit location
$ws.Range("C16").Select()
